$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update the Results column for the two new search test cases from PASS to SKIP
$ws.Range("E2").Value = "SKIP"
$ws.Range("E19").Value = "SKIP"

# Update the saved selection on the active sheet
$ws.Range("C19").Select()
